$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "ECs"
$ws.Cells.Item(2, 2).Value = "Vcam1"
$ws.Cells.Item(2, 3).Value = "Itga4"
$ws.Cells.Item(2, 4).Value = "ECs"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 1
$ws.Cells.Item(2, 7).Value = 6.664768333333334
$ws.Cells.Item(2, 8).Value = 19.994305
$ws.Cells.Item(2, 9).Value = 0.06516174319532789
$ws.Cells.Item(2, 10).Value = 0.0651617431953279
$ws.Cells.Item(2, 11).Value = 2
$ws.Cells.Item(2, 12).Value = 0.6666666666666666
$ws.Cells.Item(2, 13).Value = 0.3809099999999999
$ws.Cells.Item(2, 14).Value = 1.14273
$ws.Cells.Item(2, 15).Value = 0.006635732896411959
$ws.Cells.Item(2, 16).Value = 0.006635732896411961
$ws.Cells.Item(2, 17).Value = 2.53867690585
$ws.Cells.Item(2, 18).Value = 22.84809215265
$ws.Cells.Item(2, 19).Value = 0.0004323959229087853
$ws.Cells.Item(2, 20).Value = 0.0004323959229087856

$ws.Cells.Item(3, 1).Value = "ECs"
$ws.Cells.Item(3, 2).Value = "Vcam1"
$ws.Cells.Item(3, 3).Value = "Itga4"
$ws.Cells.Item(3, 4).Value = "FAPs"
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 1
$ws.Cells.Item(3, 7).Value = 6.664768333333334
$ws.Cells.Item(3, 8).Value = 19.994305
$ws.Cells.Item(3, 9).Value = 0.06516174319532789
$ws.Cells.Item(3, 10).Value = 0.0651617431953279
$ws.Cells.Item(3, 11).Value = 3
$ws.Cells.Item(3, 12).Value = 1
$ws.Cells.Item(3, 13).Value = 0.3194813333333333
$ws.Cells.Item(3, 14).Value = 0.9584440000000001
$ws.Cells.Item(3, 15).Value = 0.005565600255676026
$ws.Cells.Item(3, 16).Value = 0.005565600255676028
$ws.Cells.Item(3, 17).Value = 2.129269073491111
$ws.Cells.Item(3, 18).Value = 19.16342166142
$ws.Cells.Item(3, 19).Value = 0.0003626642145882124
$ws.Cells.Item(3, 20).Value = 0.0003626642145882126

$ws.Cells.Item(4, 1).Value = "ECs"
$ws.Cells.Item(4, 2).Value = "Vcam1"
$ws.Cells.Item(4, 3).Value = "Itga4"
$ws.Cells.Item(4, 4).Value = "M2"
$ws.Cells.Item(4, 5).Value = 3
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 7).Value = 6.664768333333334
$ws.Cells.Item(4, 8).Value = 19.994305
$ws.Cells.Item(4, 9).Value = 0.06516174319532789
$ws.Cells.Item(4, 10).Value = 0.0651617431953279
$ws.Cells.Item(4, 11).Value = 3
$ws.Cells.Item(4, 12).Value = 1
$ws.Cells.Item(4, 13).Value = 54.69403966666666
$ws.Cells.Item(4, 14).Value = 164.082119
$ws.Cells.Item(4, 15).Value = 0.9528104755815301
$ws.Cells.Item(4, 16).Value = 0.9528104755815303
$ws.Cells.Item(4, 17).Value = 364.5231035924772
$ws.Cells.Item(4, 18).Value = 3280.707932332295
$ws.Cells.Item(4, 19).Value = 0.06208679152366189
$ws.Cells.Item(4, 20).Value = 0.06208679152366192

$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Vcam1"
$ws.Cells.Item(5, 3).Value = "Itga4"
$ws.Cells.Item(5, 4).Value = "sCs"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 6.664768333333334
$ws.Cells.Item(5, 8).Value = 19.994305
$ws.Cells.Item(5, 9).Value = 0.06516174319532789
$ws.Cells.Item(5, 10).Value = 0.0651617431953279
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 2.008422
$ws.Cells.Item(5, 14).Value = 6.025265999999999
$ws.Cells.Item(5, 15).Value = 0.03498819126638183
$ws.Cells.Item(5, 16).Value = 0.03498819126638184
$ws.Cells.Item(5, 17).Value = 13.38566734557
$ws.Cells.Item(5, 18).Value = 120.47100611013
$ws.Cells.Item(5, 19).Value = 0.002279891534168987
$ws.Cells.Item(5, 20).Value = 0.002279891534168988

$ws.Cells.Item(6, 1).Value = "FAPs"
$ws.Cells.Item(6, 2).Value = "Vcam1"
$ws.Cells.Item(6, 3).Value = "Itga4"
$ws.Cells.Item(6, 4).Value = "ECs"
$ws.Cells.Item(6, 5).Value = 3
$ws.Cells.Item(6, 6).Value = 1
$ws.Cells.Item(6, 7).Value = 47.25592399999999
$ws.Cells.Item(6, 8).Value = 141.767772
$ws.Cells.Item(6, 9).Value = 0.4620233187619072
$ws.Cells.Item(6, 10).Value = 0.4620233187619072
$ws.Cells.Item(6, 11).Value = 2
$ws.Cells.Item(6, 12).Value = 0.6666666666666666
$ws.Cells.Item(6, 13).Value = 0.3809099999999999
$ws.Cells.Item(6, 14).Value = 1.14273
$ws.Cells.Item(6, 15).Value = 0.006635732896411959
$ws.Cells.Item(6, 16).Value = 0.006635732896411961
$ws.Cells.Item(6, 17).Value = 18.00025401083999
$ws.Cells.Item(6, 18).Value = 162.00228609756
$ws.Cells.Item(6, 19).Value = 0.003065863335217816
$ws.Cells.Item(6, 20).Value = 0.003065863335217817

$ws.Cells.Item(7, 1).Value = "FAPs"
$ws.Cells.Item(7, 2).Value = "Vcam1"
$ws.Cells.Item(7, 3).Value = "Itga4"
$ws.Cells.Item(7, 4).Value = "FAPs"
$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 1
$ws.Cells.Item(7, 7).Value = 47.25592399999999
$ws.Cells.Item(7, 8).Value = 141.767772
$ws.Cells.Item(7, 9).Value = 0.4620233187619072
$ws.Cells.Item(7, 10).Value = 0.4620233187619072
$ws.Cells.Item(7, 11).Value = 3
$ws.Cells.Item(7, 12).Value = 1
$ws.Cells.Item(7, 13).Value = 0.3194813333333333
$ws.Cells.Item(7, 14).Value = 0.9584440000000001
$ws.Cells.Item(7, 15).Value = 0.005565600255676026
$ws.Cells.Item(7, 16).Value = 0.005565600255676028
$ws.Cells.Item(7, 17).Value = 15.09738560741867
$ws.Cells.Item(7, 18).Value = 135.876470466768
$ws.Cells.Item(7, 19).Value = 0.002571437101029556
$ws.Cells.Item(7, 20).Value = 0.002571437101029558

$ws.Cells.Item(8, 1).Value = "FAPs"
$ws.Cells.Item(8, 2).Value = "Vcam1"
$ws.Cells.Item(8, 3).Value = "Itga4"
$ws.Cells.Item(8, 4).Value = "M2"
$ws.Cells.Item(8, 5).Value = 3
$ws.Cells.Item(8, 6).Value = 1
$ws.Cells.Item(8, 7).Value = 47.25592399999999
$ws.Cells.Item(8, 8).Value = 141.767772
$ws.Cells.Item(8, 9).Value = 0.4620233187619072
$ws.Cells.Item(8, 10).Value = 0.4620233187619072
$ws.Cells.Item(8, 11).Value = 3
$ws.Cells.Item(8, 12).Value = 1
$ws.Cells.Item(8, 13).Value = 54.69403966666666
$ws.Cells.Item(8, 14).Value = 164.082119
$ws.Cells.Item(8, 15).Value = 0.9528104755815301
$ws.Cells.Item(8, 16).Value = 0.9528104755815303
$ws.Cells.Item(8, 17).Value = 2584.617381740985
$ws.Cells.Item(8, 18).Value = 23261.55643566886
$ws.Cells.Item(8, 19).Value = 0.4402206580792896
$ws.Cells.Item(8, 20).Value = 0.4402206580792898

$ws.Cells.Item(9, 1).Value = "FAPs"
$ws.Cells.Item(9, 2).Value = "Vcam1"
$ws.Cells.Item(9, 3).Value = "Itga4"
$ws.Cells.Item(9, 4).Value = "sCs"
$ws.Cells.Item(9, 5).Value = 3
$ws.Cells.Item(9, 6).Value = 1
$ws.Cells.Item(9, 7).Value = 47.25592399999999
$ws.Cells.Item(9, 8).Value = 141.767772
$ws.Cells.Item(9, 9).Value = 0.4620233187619072
$ws.Cells.Item(9, 10).Value = 0.4620233187619072
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 2.008422
$ws.Cells.Item(9, 14).Value = 6.025265999999999
$ws.Cells.Item(9, 15).Value = 0.03498819126638183
$ws.Cells.Item(9, 16).Value = 0.03498819126638184
$ws.Cells.Item(9, 17).Value = 94.90983739192798
$ws.Cells.Item(9, 18).Value = 854.1885365273517
$ws.Cells.Item(9, 19).Value = 0.01616536024637011
$ws.Cells.Item(9, 20).Value = 0.01616536024637011

$ws.Cells.Item(10, 1).Value = "M2"
$ws.Cells.Item(10, 2).Value = "Vcam1"
$ws.Cells.Item(10, 3).Value = "Itga4"
$ws.Cells.Item(10, 4).Value = "ECs"
$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1
$ws.Cells.Item(10, 7).Value = 20.98736333333333
$ws.Cells.Item(10, 8).Value = 62.96209
$ws.Cells.Item(10, 9).Value = 0.2051944060881897
$ws.Cells.Item(10, 10).Value = 0.2051944060881898
$ws.Cells.Item(10, 11).Value = 2
$ws.Cells.Item(10, 12).Value = 0.6666666666666666
$ws.Cells.Item(10, 13).Value = 0.3809099999999999
$ws.Cells.Item(10, 14).Value = 1.14273
$ws.Cells.Item(10, 15).Value = 0.006635732896411959
$ws.Cells.Item(10, 16).Value = 0.006635732896411961
$ws.Cells.Item(10, 17).Value = 7.994296567299998
$ws.Cells.Item(10, 18).Value = 71.94866910569999
$ws.Cells.Item(10, 19).Value = 0.001361615270639115
$ws.Cells.Item(10, 20).Value = 0.001361615270639116

$ws.Cells.Item(11, 1).Value = "M2"
$ws.Cells.Item(11, 2).Value = "Vcam1"
$ws.Cells.Item(11, 3).Value = "Itga4"
$ws.Cells.Item(11, 4).Value = "FAPs"
$ws.Cells.Item(11, 5).Value = 3
$ws.Cells.Item(11, 6).Value = 1
$ws.Cells.Item(11, 7).Value = 20.98736333333333
$ws.Cells.Item(11, 8).Value = 62.96209
$ws.Cells.Item(11, 9).Value = 0.2051944060881897
$ws.Cells.Item(11, 10).Value = 0.2051944060881898
$ws.Cells.Item(11, 11).Value = 3
$ws.Cells.Item(11, 12).Value = 1
$ws.Cells.Item(11, 13).Value = 0.3194813333333333
$ws.Cells.Item(11, 14).Value = 0.9584440000000001
$ws.Cells.Item(11, 15).Value = 0.005565600255676026
$ws.Cells.Item(11, 16).Value = 0.005565600255676028
$ws.Cells.Item(11, 17).Value = 6.705070820884445
$ws.Cells.Item(11, 18).Value = 60.34563738796001
$ws.Cells.Item(11, 19).Value = 0.001142030038987719
$ws.Cells.Item(11, 20).Value = 0.00114203003898772

$ws.Cells.Item(12, 1).Value = "M2"
$ws.Cells.Item(12, 2).Value = "Vcam1"
$ws.Cells.Item(12, 3).Value = "Itga4"
$ws.Cells.Item(12, 4).Value = "M2"
$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 1
$ws.Cells.Item(12, 7).Value = 20.98736333333333
$ws.Cells.Item(12, 8).Value = 62.96209
$ws.Cells.Item(12, 9).Value = 0.2051944060881897
$ws.Cells.Item(12, 10).Value = 0.2051944060881898
$ws.Cells.Item(12, 11).Value = 3
$ws.Cells.Item(12, 12).Value = 1
$ws.Cells.Item(12, 13).Value = 54.69403966666666
$ws.Cells.Item(12, 14).Value = 164.082119
$ws.Cells.Item(12, 15).Value = 0.9528104755815301
$ws.Cells.Item(12, 16).Value = 0.9528104755815303
$ws.Cells.Item(12, 17).Value = 1147.883682652079
$ws.Cells.Item(12, 18).Value = 10330.95314386871
$ws.Cells.Item(12, 19).Value = 0.1955113796515577
$ws.Cells.Item(12, 20).Value = 0.1955113796515577

$ws.Cells.Item(13, 1).Value = "M2"
$ws.Cells.Item(13, 2).Value = "Vcam1"
$ws.Cells.Item(13, 3).Value = "Itga4"
$ws.Cells.Item(13, 4).Value = "sCs"
$ws.Cells.Item(13, 5).Value = 3
$ws.Cells.Item(13, 6).Value = 1
$ws.Cells.Item(13, 7).Value = 20.98736333333333
$ws.Cells.Item(13, 8).Value = 62.96209
$ws.Cells.Item(13, 9).Value = 0.2051944060881897
$ws.Cells.Item(13, 10).Value = 0.2051944060881898
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 2.008422
$ws.Cells.Item(13, 14).Value = 6.025265999999999
$ws.Cells.Item(13, 15).Value = 0.03498819126638183
$ws.Cells.Item(13, 16).Value = 0.03498819126638184
$ws.Cells.Item(13, 17).Value = 42.15148224066
$ws.Cells.Item(13, 18).Value = 379.36334016594
$ws.Cells.Item(13, 19).Value = 0.007179381127005206
$ws.Cells.Item(13, 20).Value = 0.007179381127005209

$ws.Cells.Item(14, 1).Value = "sCs"
$ws.Cells.Item(14, 2).Value = "Vcam1"
$ws.Cells.Item(14, 3).Value = "Itga4"
$ws.Cells.Item(14, 4).Value = "ECs"
$ws.Cells.Item(14, 5).Value = 3
$ws.Cells.Item(14, 6).Value = 1
$ws.Cells.Item(14, 7).Value = 27.37233166666667
$ws.Cells.Item(14, 8).Value = 82.116995
$ws.Cells.Item(14, 9).Value = 0.2676205319545753
$ws.Cells.Item(14, 10).Value = 0.2676205319545753
$ws.Cells.Item(14, 11).Value = 2
$ws.Cells.Item(14, 12).Value = 0.6666666666666666
$ws.Cells.Item(14, 13).Value = 0.3809099999999999
$ws.Cells.Item(14, 14).Value = 1.14273
$ws.Cells.Item(14, 15).Value = 0.006635732896411959
$ws.Cells.Item(14, 16).Value = 0.006635732896411961
$ws.Cells.Item(14, 17).Value = 10.42639485515
$ws.Cells.Item(14, 18).Value = 93.83755369634999
$ws.Cells.Item(14, 19).Value = 0.001775858367646243
$ws.Cells.Item(14, 20).Value = 0.001775858367646244

$ws.Cells.Item(15, 1).Value = "sCs"
$ws.Cells.Item(15, 2).Value = "Vcam1"
$ws.Cells.Item(15, 3).Value = "Itga4"
$ws.Cells.Item(15, 4).Value = "FAPs"
$ws.Cells.Item(15, 5).Value = 3
$ws.Cells.Item(15, 6).Value = 1
$ws.Cells.Item(15, 7).Value = 27.37233166666667
$ws.Cells.Item(15, 8).Value = 82.116995
$ws.Cells.Item(15, 9).Value = 0.2676205319545753
$ws.Cells.Item(15, 10).Value = 0.2676205319545753
$ws.Cells.Item(15, 11).Value = 3
$ws.Cells.Item(15, 12).Value = 1
$ws.Cells.Item(15, 13).Value = 0.3194813333333333
$ws.Cells.Item(15, 14).Value = 0.9584440000000001
$ws.Cells.Item(15, 15).Value = 0.005565600255676026
$ws.Cells.Item(15, 16).Value = 0.005565600255676028
$ws.Cells.Item(15, 17).Value = 8.744949017308889
$ws.Cells.Item(15, 18).Value = 78.70454115578001
$ws.Cells.Item(15, 19).Value = 0.001489468901070538
$ws.Cells.Item(15, 20).Value = 0.001489468901070539

$ws.Cells.Item(16, 1).Value = "sCs"
$ws.Cells.Item(16, 2).Value = "Vcam1"
$ws.Cells.Item(16, 3).Value = "Itga4"
$ws.Cells.Item(16, 4).Value = "M2"
$ws.Cells.Item(16, 5).Value = 3
$ws.Cells.Item(16, 6).Value = 1
$ws.Cells.Item(16, 7).Value = 27.37233166666667
$ws.Cells.Item(16, 8).Value = 82.116995
$ws.Cells.Item(16, 9).Value = 0.2676205319545753
$ws.Cells.Item(16, 10).Value = 0.2676205319545753
$ws.Cells.Item(16, 11).Value = 3
$ws.Cells.Item(16, 12).Value = 1
$ws.Cells.Item(16, 13).Value = 54.69403966666666
$ws.Cells.Item(16, 14).Value = 164.082119
$ws.Cells.Item(16, 15).Value = 0.9528104755815301
$ws.Cells.Item(16, 16).Value = 0.9528104755815303
$ws.Cells.Item(16, 17).Value = 1497.103393945823
$ws.Cells.Item(16, 18).Value = 13473.9305455124
$ws.Cells.Item(16, 19).Value = 0.254991646327021
$ws.Cells.Item(16, 20).Value = 0.254991646327021

$ws.Cells.Item(17, 1).Value = "sCs"
$ws.Cells.Item(17, 2).Value = "Vcam1"
$ws.Cells.Item(17, 3).Value = "Itga4"
$ws.Cells.Item(17, 4).Value = "sCs"
$ws.Cells.Item(17, 5).Value = 3
$ws.Cells.Item(17, 6).Value = 1
$ws.Cells.Item(17, 7).Value = 27.37233166666667
$ws.Cells.Item(17, 8).Value = 82.116995
$ws.Cells.Item(17, 9).Value = 0.2676205319545753
$ws.Cells.Item(17, 10).Value = 0.2676205319545753
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 2.008422
$ws.Cells.Item(17, 14).Value = 6.025265999999999
$ws.Cells.Item(17, 15).Value = 0.03498819126638183
$ws.Cells.Item(17, 16).Value = 0.03498819126638184
$ws.Cells.Item(17, 17).Value = 54.97519311063
$ws.Cells.Item(17, 18).Value = 494.77673799567
$ws.Cells.Item(17, 19).Value = 0.00936355835883753
$ws.Cells.Item(17, 20).Value = 0.009363558358837532

